$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 116, shifting the existing Cereza records
# (old rows 116-122) down to 117-123, then populate the new row 116 with
# the "Early Burlat" entry.
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value = 44516
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100103
$ws.Range("H116").Value = "Frutos de hueso (carozo)"
$ws.Range("I116").Value = 100103001
$ws.Range("J116").Value = "Cereza"
$ws.Range("K116").Value = "Early Burlat"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 25
$ws.Range("N116").Value = 3500
$ws.Range("O116").Value = 3500
$ws.Range("P116").Value = 3500
$ws.Range("Q116").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R116").Value = "Región del Maule"
$ws.Range("S116").Value = 3500
$ws.Range("T116").Value = 1
